{"js": "// Replace the 25 division-problem answer strings in the single table of\n// the document, in document order (row-major, left-to-right), leaving the\n// title paragraph and the empty \"scratch work\" rows untouched.\nconst newValues = [\n  \"97\u00f72=48, 1\",\n  \"33\u00f74=8, 1\",\n  \"89\u00f72=44, 1\",\n  \"75\u00f78=9, 3\",\n  \"55\u00f79=6, 1\",\n  \"89\u00f76=14, 5\",\n  \"66\u00f73=22, 0\",\n  \"42\u00f79=4, 6\",\n  \"70\u00f74=17, 2\",\n  \"47\u00f77=6, 5\",\n  \"70\u00f79=7, 7\",\n  \"91\u00f73=30, 1\",\n  \"25\u00f74=6, 1\",\n  \"35\u00f78=4, 3\",\n  \"17\u00f76=2, 5\",\n  \"37\u00f74=9, 1\",\n  \"99\u00f74=24, 3\",\n  \"34\u00f74=8, 2\",\n  \"35\u00f75=7, 0\",\n  \"19\u00f76=3, 1\",\n  \"63\u00f78=7, 7\",\n  \"53\u00f76=8, 5\",\n  \"89\u00f76=14, 5\",\n  \"56\u00f76=9, 2\",\n  \"46\u00f76=7, 4\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather every cell's first paragraph up front (load text so we know which\n// cells actually hold an answer string vs. the blank scratch-work cells).\nconst cellParas = [];\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < cells.items.length; c++) {\n    const paras = cells.items[c].body.paragraphs;\n    paras.load(\"items\");\n    cellParas.push(paras);\n  }\n}\nawait context.sync();\n\nconst firstParas = cellParas.map((paras) => paras.items[0]);\nfirstParas.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet idx = 0;\nfor (const p of firstParas) {\n  if (p.text !== \"\") {\n    p.insertText(newValues[idx], \"Replace\");\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answer strings in the single table of\n# the document, in document order (row-major, left-to-right), leaving the\n# title paragraph and the empty \"scratch work\" rows untouched.\n$newValues = @(\n  \"97\u00f72=48, 1\",\n  \"33\u00f74=8, 1\",\n  \"89\u00f72=44, 1\",\n  \"75\u00f78=9, 3\",\n  \"55\u00f79=6, 1\",\n  \"89\u00f76=14, 5\",\n  \"66\u00f73=22, 0\",\n  \"42\u00f79=4, 6\",\n  \"70\u00f74=17, 2\",\n  \"47\u00f77=6, 5\",\n  \"70\u00f79=7, 7\",\n  \"91\u00f73=30, 1\",\n  \"25\u00f74=6, 1\",\n  \"35\u00f78=4, 3\",\n  \"17\u00f76=2, 5\",\n  \"37\u00f74=9, 1\",\n  \"99\u00f74=24, 3\",\n  \"34\u00f74=8, 2\",\n  \"35\u00f75=7, 0\",\n  \"19\u00f76=3, 1\",\n  \"63\u00f78=7, 7\",\n  \"53\u00f76=8, 5\",\n  \"89\u00f76=14, 5\",\n  \"56\u00f76=9, 2\",\n  \"46\u00f76=7, 4\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n# First pass: snapshot which cells currently hold an answer string (vs. the\n# blank scratch-work cells), in row-major document order, so edits made\n# along the way never influence which cells we still treat as \"to do\".\n$targets = @()\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $txt = $cell.Range.Text\n    # A Word table cell's Range.Text always ends with the cell-mark (\\r\\a);\n    # strip those two trailing characters before checking/comparing content.\n    $clean = $txt.Substring(0, $txt.Length - 2)\n    if ($clean -ne \"\") {\n      $targets += $cell\n    }\n  }\n}\n\nfor ($i = 0; $i -lt $targets.Count; $i++) {\n  $targets[$i].Range.Text = $newValues[$i]\n}\n"}
